$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Professional summary paragraph: neutralize the "Black and
#    Asian-American voters" phrasing to "50M voters".
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters, developed",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2. Rebuild the "KEY ACHIEVEMENTS AND IMPACT" bullet list.
#    Locate the four existing bullet paragraphs that sit right after the
#    "Impact" sub-heading, just before "TECHNICAL SKILLS".
# ---------------------------------------------------------------------
$boldColor = 5258796   # RGB(0x2C,0x3E,0x50) packed as BGR for Font.Color

function Set-BoldRun($paragraphIndex, $searchText) {
    $p = $d.Paragraphs($paragraphIndex)
    $rng = $d.Range($p.Range.Start, $p.Range.End)
    $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "WARNING: could not find '$searchText' in paragraph $paragraphIndex"
    }
    $rng.Font.Bold = 1
    $rng.Font.Color = $boldColor
}

# Find the anchor paragraphs by content so this is resilient to any
# paragraph-number drift caused by step 1. Paragraph.Range.Text always
# carries a trailing "\r" paragraph mark, so trim it before comparing.
$implIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text.TrimEnd("`r")
    if ($ptext -eq "Impact") {
        $implIdx = $i
    }
}

$p1 = $d.Paragraphs($implIdx + 1)   # "Discovered systematic race coding errors..."
$p2 = $d.Paragraphs($implIdx + 2)   # "Algorithm reduced mapping costs by 73.5%..."
$p3 = $d.Paragraphs($implIdx + 3)   # "Built redistricting platform..."
$p4 = $d.Paragraphs($implIdx + 4)   # "Achieved 87% ... 71%"

# --- Remove paragraph 1 ("Discovered systematic race coding errors...") ---
$d.Range($p1.Range.Start, $p1.Range.End).Delete() | Out-Null

# --- Paragraph 2 becomes: "...mapping costs **73.5%**" and split off the
#     $4.7M clause into its own new paragraph. ---
$p2 = $d.Paragraphs($implIdx + 1)
$p2.Range.Text = "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs **73.5%**"
Set-BoldRun ($implIdx + 1) "73.5%"

$p2 = $d.Paragraphs($implIdx + 1)
$p2.Range.InsertParagraphAfter() | Out-Null
$pNew = $d.Paragraphs($implIdx + 2)
$pNew.Range.Text = "• **`$4.7M** savings enabled nonprofit access"
Set-BoldRun ($implIdx + 2) "`$4.7M"

# --- Paragraph 3 ("Built redistricting platform...") becomes the Supreme
#     Court legal-precedent bullet, with two extra bullets inserted after
#     it. ---
$p3 = $d.Paragraphs($implIdx + 3)
$p3.Range.Text = "• Legal precedent: Data analysis utilized in Supreme Court case"

$p3 = $d.Paragraphs($implIdx + 3)
$p3.Range.InsertParagraphAfter() | Out-Null
$pNew = $d.Paragraphs($implIdx + 4)
$pNew.Range.Text = "• Expert methodology validated at highest judicial level"

$pNew.Range.InsertParagraphAfter() | Out-Null
$pNew2 = $d.Paragraphs($implIdx + 5)
$pNew2.Range.Text = "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"

# --- Paragraph 4 ("Achieved 87% ... 71%") becomes the racial
#     classification accuracy-improvement bullet. ---
$p4 = $d.Paragraphs($implIdx + 6)
$p4.Range.Text = "• **178%** accuracy improvement in racial classification algorithms"
Set-BoldRun ($implIdx + 6) "178%"

Write-Output "done"
